$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing old rows 13-23 down to 14-24.
# (A new "Docentes responsaveis" value row is being introduced under that label.)
$ws.Rows.Item(13).Insert()

# The insert copied formatting from the row above into A13; the target layout
# has no cell in column A on this row, so clear it back out.
$ws.Range("A13").ClearContents()

# Row 10 (Objetivos:) previously (incorrectly) held the teacher name in B/C;
# put the real course-objectives paragraph there instead.
$ws.Range("B10:C10").Value = "Fornecer ao aluno os conhecimentos fundamentais relativos ao cultivo de células animais, visando seu emprego como instrumento de obtenção de produtos biotecnológicos de alto valor agregado"

# Row 13 (Docentes responsaveis:) now gets the teacher name that used to sit
# (mistakenly) under Objetivos.
$ws.Range("B13:C13").Value = "101761 - Arnaldo Márcio Ramalho Prata"

# Row 14 (Programa resumido:) previously held the placeholder "Semestral";
# replace it with the actual short-syllabus summary (Portuguese).
$ws.Range("B14:C14").Value = "Introdução à Tecnologia de Cultivo de Células Animais, Crescimento e Morte de Células Animais Cultivadas in vitro, Biorreatores para Células Animais, Aplicações do Cultivo de Células Animais"

# Row 16 (Programa:) previously held a stray date; replace it with the full
# numbered Portuguese syllabus text.
$ws.Range("B16:C16").Value = "1. Introdução à Tecnologia de Cultivo de Células Animais  Principais marcos e razões da cultura de células animais, Tipos de culturas de células animais, Emprego de células animais. 2. Mecanismo de Crescimento e Morte de Células Animais Cultivadas in vitro  Mecanismos de proliferação celular, Mecanismos de morte celular, Influência das condições ambientais sobre a morte celular, Métodos de detecção da morte celular, Controle da apoptose por técnicas moleculares. 3. Biorreatores para Células Animais  Propagação de inóculo e sistemas de cultivo em pequena escala, Tipos de biorreatores, Aeração e agitação, Aspectos econômicos na seleção de biorreatores. 4. Aplicações do Cultivo de Células Animais  Proteinas recombinantes terapêuticas, Anticorpos monoclonais, Vacinas virais, Bioinseticidas, Terapias celulares e células-tronco, Terapia gênica"

# Row 19 (Metodo:) previously (incorrectly) held the teacher name; replace it
# with the evaluation method text.
$ws.Range("B19:C19").Value = "A avaliação será feita por meio de provas escritas."

# Row 20 (Criterio:) now holds the final-grade formula text.
$ws.Range("B20:C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = [P1 +(2 x P2)] / 3"

# Row 21 (Norma de recuperacao:) now holds the recovery-grade formula text.
$ws.Range("B21:C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"

# Row 22 (Bibliografia:) previously held the recovery-grade text by mistake;
# replace it with the actual bibliography list.
$ws.Range("B22:C22").Value = "1. MORAES, A.M., AUGUSTO, E.F.P., CASTILHO, L.R. Tecnologia do Cultivo de Células Animais – de Biofármacos a Terapia Gênica. São Paulo: Rocca, 2008.2. VITOLO, M. (Coordenador). Biotecnologia Farmacêutica – Aspectos sobre aplicação industrial. São Paulo: Edgard Blücher Ltda, 2015.3. SHULER, M.L., KARGI, F. Bioprocess Engineering – Basic Concepts. Second edition. New Jersey: Prentice Hall, 2002."

# Column A no longer shares a combined width definition with column B; give it
# its own explicit width matching the original (column B keeps its own too).
$ws.Columns.Item(1).ColumnWidth = 30.7109375
$ws.Columns.Item(2).ColumnWidth = 60.7109375
